$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update fitted-distribution values (final fitted distributions 3/12)
$ws.Range("B2").Value = 0.83
$ws.Range("C2").Value = 0.77
$ws.Range("D2").Value = 0.89

$ws.Range("B3").Value = 0.67
$ws.Range("C3").Value = 0.56
$ws.Range("D3").Value = 0.74

$ws.Range("C4").Value = 0.74
$ws.Range("D4").Value = 0.94

$ws.Range("D5").Value = 0.99

# Update the active selection to match the saved sheet view state
$ws.Range("E12").Select()
